$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 3, pushing the existing row 3
# (date 44636, Primera quality data) down to row 5.
$ws.Rows.Item(3).Resize(2).Insert()

# Row 3: updated data for the same report date group (new sample, Primera quality)
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(3, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(3, 4).Value = 44650
$ws.Cells.Item(3, 4).NumberFormat = $ws.Cells.Item(2, 4).NumberFormat
$ws.Cells.Item(3, 5).Value = 15
$ws.Cells.Item(3, 6).Value = "Fruta"
$ws.Cells.Item(3, 7).Value = 100107
$ws.Cells.Item(3, 8).Value = "Otros"
$ws.Cells.Item(3, 9).Value = 100107011
$ws.Cells.Item(3, 10).Value = "Tuna"
$ws.Cells.Item(3, 11).Value = "Sin especificar"
$ws.Cells.Item(3, 12).Value = "Primera"
$ws.Cells.Item(3, 13).Value = 160
$ws.Cells.Item(3, 14).Value = 31000
$ws.Cells.Item(3, 15).Value = 32000
$ws.Cells.Item(3, 16).Value = 31500
$ws.Cells.Item(3, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(3, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(3, 19).Value = 1575
$ws.Cells.Item(3, 20).Value = 20

# Row 4: new row, Segunda quality
$ws.Cells.Item(4, 1).Value = 1
$ws.Cells.Item(4, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(4, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(4, 4).Value = 44650
$ws.Cells.Item(4, 4).NumberFormat = $ws.Cells.Item(2, 4).NumberFormat
$ws.Cells.Item(4, 5).Value = 15
$ws.Cells.Item(4, 6).Value = "Fruta"
$ws.Cells.Item(4, 7).Value = 100107
$ws.Cells.Item(4, 8).Value = "Otros"
$ws.Cells.Item(4, 9).Value = 100107011
$ws.Cells.Item(4, 10).Value = "Tuna"
$ws.Cells.Item(4, 11).Value = "Sin especificar"
$ws.Cells.Item(4, 12).Value = "Segunda"
$ws.Cells.Item(4, 13).Value = 250
$ws.Cells.Item(4, 14).Value = 29000
$ws.Cells.Item(4, 15).Value = 30000
$ws.Cells.Item(4, 16).Value = 29500
$ws.Cells.Item(4, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(4, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(4, 19).Value = 1475
$ws.Cells.Item(4, 20).Value = 20
